$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.937.40"
$ws.Range("E2").Value = "  -0.01%  "

$ws.Range("D3").Value = "1.743.63"
$ws.Range("E3").Value = "  -0.74%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.29%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.24%  "

$ws.Range("E6").Value = "  -0.21%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5242"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.80%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2752"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.69%  "

$ws.Range("E9").Value = "  -2.40%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06154"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.20%  "

$ws.Range("D11").Value = "1.740.99"
$ws.Range("E11").Value = "  -1.70%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07098"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.93%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.21"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.04%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6435"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.68%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.527"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.21%  "

$ws.Range("E16").Value = "  -0.48%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9997"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.27%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.000"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.19%  "

$ws.Range("D19").Value = "25.913.37"
$ws.Range("E19").Value = "  -0.16%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.54"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.34%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000006679"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.13%  "

$ws.Range("D22").Value = "1.961.23"
$ws.Range("E22").Value = "  -1.78%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.303"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.07%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.767"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.49%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.173"
$ws.Range("D25").Style = "Normal"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "140.68"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.28%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.521"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.01%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.19"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.65%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.801"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.03%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "102.76"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.36%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08340"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.37%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.732"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.85%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.562"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.23%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04539"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.84%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.617"
$ws.Range("D35").Style = "Normal"

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9804"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.20%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6203"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.81%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.690"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.91%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01590"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.07%  "

$ws.Range("E40").Value = "  -1.34%  "

$ws.Range("E41").Value = "  -0.20%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "100.29"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.13%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.3876"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.53%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.7343"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.96%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.018"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.20%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.05336"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.14%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1126"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.32%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.276"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.14%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "53.63"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.55%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "30.21"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.48%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.704"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.21%  "
